$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (H1, style index 1: bold/centered/bordered) onto the
# two new header cells before writing their text, so I1/J1 end up sharing
# the same cellXf as the existing headers (B1:H1) instead of minting a new one.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
